$d = $word.ActiveDocument

# Run-properties block shared by every run touched by this edit (same
# Lucida Sans / 24-half-point formatting the placeholder runs already had).
$rPr = '<w:rPr><w:rFonts w:ascii="Lucida Sans" w:hAnsi="Lucida Sans"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

function Build-RunXml($text, [bool]$preserve) {
    $escaped = $text.Replace('&', '&amp;').Replace('<', '&lt;').Replace('>', '&gt;')
    $spaceAttr = ''
    if ($preserve) { $spaceAttr = ' xml:space="preserve"' }
    return '<w:r>' + $rPr + '<w:t' + $spaceAttr + '>' + $escaped + '</w:t></w:r>'
}

function Replace-Placeholder($paragraphIndex, $segments) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $rng = $p.Range

    if ($rng.Text -ne "Ipsum`r") {
        throw "Paragraph $paragraphIndex was expected to contain the 'Ipsum' placeholder but contains '$($rng.Text)'."
    }

    $owx = $rng.WordOpenXML
    $openTag = [regex]::Match($owx, '<w:p\s[^>]*>').Value
    $pPr = [regex]::Match($owx, '<w:pPr>.*?</w:pPr>').Value

    $runsXml = ''
    foreach ($seg in $segments) {
        $runsXml += Build-RunXml $seg[0] $seg[1]
    }

    $openTagNs = $openTag.Substring(0, 4) + ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"' + $openTag.Substring(4)

    $xml = $openTagNs + $pPr + $runsXml + '</w:p>'
    $rng.InsertXML($xml) | Out-Null
}

# --- Paragraph: "Model View Controller (MVC)" justification -------------
$mvcSegments = @(
    ,@("The Model View Controller pattern was chosen due to its attribute of introducing the concept of a ‘controller’. Since ", $true)
    ,@("we will have to", $false)
    ,@(" be using multiple flows of data, separating the user’s action of requesting the data and presenting the data will make the implementation easier. ", $true)
    ,@("The controller will assist in this manner by making the system centralized.  ", $true)
    ,@("Distinguishing the separation will allow it ", $true)
    ,@("to be", $false)
    ,@(" more efficient at processing the data before sending it to be viewed", $true)
    ,@(", as well as assuring a seamless experience. ", $true)
)
Replace-Placeholder 9 $mvcSegments

# --- Paragraph: Class diagram justification ------------------------------
$classSegments = @(
    ,@("Class diagram was chosen since the concept of a ‘Course’ in our system resembled a class the most. ", $true)
    ,@("Since the ‘Course’ itself would have attributes like a class, it made sense to use the Class diagram for our detail design. ", $true)
)
Replace-Placeholder 11 $classSegments
